# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" contains the daily conversion note in A1
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 6.99 = 27916.08 pesos`n✅ 27916.08 pesos = 6.97 = 968.18 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# Sheet "tasas" contains the updated rate figures
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 143
$wsTasas.Range("O10").Value = 3992
$wsTasas.Range("N12").Value = 4005
$wsTasas.Range("O12").Value = 138.9
